# Horarios actualizados Línea 141 - 1146
# Update scrape timestamp ("Última actualización"), refresh the arrivals
# table on each sheet, and append a new arrival row on the main "LP1912"
# sheet.

$wb = $excel.ActiveWorkbook

$newScrapTime = "03:08:51"

# ---------------------------------------------------------------------
# Sheet "LP1912" (main schedule) - rows shift up one arrival, and a new
# arrival ("11_ETCHEVERRY") is appended as row 9.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newScrapTime"
$ws1.Range("A3").Value = "Total filas: 4"

$ws1.Range("A6").Value = $newScrapTime
$ws1.Range("B6").Value = "04:01"
$ws1.Range("C6").Value = "81_EL PELIGRO"
$ws1.Range("D6").Value = 53
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = $newScrapTime
$ws1.Range("B7").Value = "04:45"
$ws1.Range("C7").Value = "215A_EL PATO"
$ws1.Range("D7").Value = 97
$ws1.Range("E7").Value = "LP1912"

$ws1.Range("A8").Value = $newScrapTime
$ws1.Range("B8").Value = "04:48"
$ws1.Range("C8").Value = "14_ABASTO"
$ws1.Range("D8").Value = 100
$ws1.Range("E8").Value = "LP1912"

$ws1.Range("A9").Value = $newScrapTime
$ws1.Range("B9").Value = "04:53"
$ws1.Range("C9").Value = "11_ETCHEVERRY"
$ws1.Range("D9").Value = 105
$ws1.Range("E9").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "LP1912-215" - single arrival refreshed in place.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newScrapTime"

$ws2.Range("A6").Value = $newScrapTime
$ws2.Range("B6").Value = "04:45"
$ws2.Range("C6").Value = "215A_EL PATO"
$ws2.Range("D6").Value = 97
$ws2.Range("E6").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173" - no arrivals, only the scrape timestamp refreshes.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newScrapTime"
